# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# "ff84103e-e205-4350-9f14-47d23e1a5a60.md" file row (row 3 in each table)
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-07 13:25:24"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-07 13:24:59"
$wsZhCn.Range("K3").Value = "2016-09-07 13:25:54"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-07 13:25:24"
$wsDeDe.Range("K3").Value = "2016-09-07 13:26:24"
